$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.306.20"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.538.80"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "3.536.40"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -4.09%  "
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "4.136.09"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  -4.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("D16").Value = "3.540.88"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "66.381.61"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").Value = "3.675.43"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "3.527.21"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E37").Value = "  -3.10%  "
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0858"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -7.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("E51").Value = "  -5.43%  "
